# Update "想去人数" (want-to-go count) figures for several expo rows.
# Applies to sheet "展览" (Exhibition) and sheet "全部类型" (All types).

$wb = $excel.ActiveWorkbook

# -- Sheet "展览": rows 3-7 in column F --
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 89
$wsExhibit.Range("F4").Value = 46
$wsExhibit.Range("F5").Value = 2382
$wsExhibit.Range("F6").Value = 221
$wsExhibit.Range("F7").Value = 374

# -- Sheet "全部类型": rows 3,4,5,6 and 9 in column F --
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 89
$wsAll.Range("F4").Value = 46
$wsAll.Range("F5").Value = 2382
$wsAll.Range("F6").Value = 221
$wsAll.Range("F9").Value = 374
